$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("B1").Value = "Groepnaam"

# --- Row 2 (style stays like row1/row3, value is a pure-digit code so must
#     land as text without disturbing the cell's number format/style) ---
$ws.Range("A2").Formula = "=""460003000"""
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("B2").Value = "Brandblusmiddelen"

# --- Row 3 ---
$ws.Range("A3").Value = "460A16100"
$ws.Range("B3").Value = "Kopschot wandstelling"

# --- Row 4 ---
$ws.Range("A4").Value = "460A16000"
$ws.Range("B4").Value = "Kopschot wandstelling"

# --- Row 5: needs the "plain left" style (same as the predefined-but-unused
#     style already in styles.xml) instead of the wrap/vcenter one used by
#     rows 1-4. Build that exact alignment on a scratch cell first and copy
#     just the formatting across so the engine reuses the existing style
#     slot instead of minting a new one. ---
$ws.Range("Z100").HorizontalAlignment = -4131
$ws.Range("Z100").VerticalAlignment = -4107
$ws.Range("Z100").WrapText = $false
$ws.Range("Z100").Copy()
$ws.Range("A5:A9").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

$ws.Range("A5").Formula = "=""460000300"""
$ws.Range("A5").Copy()
$ws.Range("A5").PasteSpecial(-4163)
$ws.Range("B5").Value = "Belijning vloer"

# --- Row 6 ---
$ws.Range("A6").Value = "130B01600"
$ws.Range("B6").Value = "Platenwagen t.b.v. bake-off"

# --- Row 7 ---
$ws.Range("A7").Value = "130B50001"
$ws.Range("B7").Value = "Brood werktafel"

# --- Row 8 ---
$ws.Range("A8").Value = "560B01000"
$ws.Range("B8").Value = "Winkelwagens ELA"

# --- Row 9 ---
$ws.Range("A9").Value = "560B01300"
$ws.Range("B9").Value = "Winkelwagens Tango"

# --- Column B width (Excel "best fit" width after the new text was added) ---
$ws.Columns.Item(2).ColumnWidth = 23.1015625

# --- Selection, matches the saved file ---
[void]$ws.Range("B4").Select()
